# Update odds values on Sheet1 for rows 3, 4, and 5 to match the
# latest FlashScore snapshot (commit: "Atualizando o arquivo XLSX").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("G3").Value = 1.48
$ws.Range("H3").Value = 4.2
$ws.Range("I3").Value = 6.5
$ws.Range("J3").Value = 2.05
$ws.Range("K3").Value = 2.3
$ws.Range("M3").Value = 1.04
$ws.Range("N3").Value = 13
$ws.Range("O3").Value = 1.25
$ws.Range("P3").Value = 3.75
$ws.Range("Q3").Value = 1.88
$ws.Range("R3").Value = 1.98
$ws.Range("S3").Value = 1.36
$ws.Range("T3").Value = 3
$ws.Range("AD3").Value = 8
$ws.Range("AE3").Value = 21
$ws.Range("AT3").Value = 3
$ws.Range("BC3").Value = 301

# Row 4
$ws.Range("G4").Value = 2.75
$ws.Range("I4").Value = 2.15
$ws.Range("J4").Value = 3.2
$ws.Range("L4").Value = 2.67
$ws.Range("O4").Value = 1.18
$ws.Range("P4").Value = 4.31
$ws.Range("W4").Value = 10
$ws.Range("X4").Value = 13.5
$ws.Range("Y4").Value = 8.75
$ws.Range("Z4").Value = 26
$ws.Range("AA4").Value = 17
$ws.Range("AB4").Value = 19.5
$ws.Range("AH4").Value = 8.5
$ws.Range("AI4").Value = 10.25
$ws.Range("AK4").Value = 17.5
$ws.Range("AL4").Value = 13
$ws.Range("AN4").Value = 5
$ws.Range("AS4").Value = 175
$ws.Range("AX4").Value = 4.35

# Row 5
$ws.Range("N5").Value = 6.95
